# REVER_DailyTracker_NIRMAL.xlsx - "Add files via upload"
#
# NOV-2020 sheet (A1:G23) gets 4 new daily-log rows (8-11, covering
# 07-Nov-2020 .. 10-Nov-2020) appended after the existing 6 rows of data,
# plus a brand-new "Comments" shared string used by row 10.
#
# Rows 8 & 9 are "Week off" rows (same look as row 2).
# Row 10 is a normal completed-task row (same look as rows 3-7) whose task
#   text is a brand-new shared string.
# Row 11 only has the No./Date filled in (task not started yet), matching
#   the blank-D-cell look used in rows 3-7 (just without any value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2020")

# ---------------------------------------------------------------------
# Row 8: No=7, Date=07-Nov-2020, Task="Week off"  (styled like row 2)
# ---------------------------------------------------------------------
$ws.Range("A2:G2").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G2").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("F8").PasteSpecial(-4122)

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 44142
$ws.Cells.Item(8, 4).Value = "Week off"

# ---------------------------------------------------------------------
# Row 9: No=8, Date=08-Nov-2020, Task="Week off"  (styled like row 2)
# ---------------------------------------------------------------------
$ws.Range("A2:G2").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("F9").PasteSpecial(-4122)

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 44143
$ws.Cells.Item(9, 4).Value = "Week off"

# ---------------------------------------------------------------------
# Row 10: No=9, Date=09-Nov-2020, completed task (styled like row 7)
# ---------------------------------------------------------------------
$ws.Range("A7:G7").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 44144
$ws.Cells.Item(10, 4).Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing, Retesting on B2C/B2B app and Sonia application"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = "Completed"
$ws.Rows.Item(10).RowHeight = 30

# ---------------------------------------------------------------------
# Row 11: No=10, Date=10-Nov-2020, task not started yet (D left blank)
# ---------------------------------------------------------------------
$ws.Range("A7:G7").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Cells.Item(11, 4).ClearContents()

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 44145

# Selection moved down one row as the user kept editing the sheet.
$ws.Range("D8").Select() | Out-Null
